$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the day label and the section label text
$ws.Range("C5").Value = "Tuesday"
$ws.Range("C7").Value = "asdsads"

# Remove the now-unused "Bench Press" workout row (row 9) and the blank
# spacer row above it (row 8).
$ws.Range("A8:A9").EntireRow.Delete()

# Columns D and E no longer hold any data; resize column C to fit the
# shorter remaining content.
$ws.Columns.Item(3).ColumnWidth = 7.2
